# Insert a new data row at row 161 (pushes existing rows 161-404 down to 162-405)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(161).Insert()

$ws.Range("A161").Value = 3
$ws.Range("B161").Value = "Femacal de La Calera"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 45036
$ws.Range("E161").Value = 5
$ws.Range("F161").Value = 100112039
$ws.Range("G161").Value = "Ciboulette"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 120
$ws.Range("K161").Value = 1500
$ws.Range("L161").Value = 1500
$ws.Range("M161").Value = 1500
$ws.Range("N161").Value = "$/docena de atados"
$ws.Range("O161").Value = "Provincia de Quillota"
$ws.Range("P161").Value = 500
$ws.Range("Q161").Value = 3
$ws.Range("R161").Value = "Hortaliza"
